# This edit re-shuffles the weekly price records: for a set of rows, the
# values in columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio
# maximo), P (Precio promedio ponderado) and S (Precio $/Kg) are replaced
# by the values that, before the edit, belonged to another row in the same
# column set (a permutation of rows). Rows 3 and 10 are left untouched.
#
# Because several rows swap/rotate values with each other, we must first
# snapshot every relevant cell's original value before writing any new
# value, otherwise a write to one row could clobber data still needed as
# the source for another row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "M", "N", "O", "P", "S")

# Maps each destination row to the row whose original (pre-edit) values it
# should receive.
$rowMap = @{
    2  = 32
    4  = 41
    5  = 13
    6  = 24
    7  = 23
    8  = 31
    9  = 21
    11 = 26
    12 = 30
    13 = 5
    14 = 19
    15 = 38
    16 = 9
    17 = 16
    18 = 28
    19 = 34
    20 = 33
    21 = 27
    22 = 14
    23 = 35
    24 = 29
    25 = 40
    26 = 2
    27 = 11
    28 = 7
    29 = 17
    30 = 22
    31 = 12
    32 = 18
    33 = 20
    34 = 15
    35 = 37
    36 = 6
    37 = 36
    38 = 25
    39 = 4
    40 = 8
    41 = 39
}

# 1) Snapshot original values for every row that participates in the
#    permutation (both as source and destination -- here that's the same
#    set of rows).
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowValues = @{}
    foreach ($c in $cols) {
        $rowValues[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowValues
}

# 2) Apply the permutation using the snapshot as the source of truth.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcValues[$c]
    }
}
